# Apply the attendance_reports sync edit to the "Session Analysis Results" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Class Statistics: Missing / Pending session counts (L7, L8) ---
$ws.Range("L7").Value = 33
$ws.Range("L8").Value = 60

# --- 2. "Recorded By" column: swap the order of "dnasr281@gmail.com, System" ---
#         to "System, dnasr281@gmail.com" for every affected session row.
$swappedRecordedByRows = @(8,9,10,12,14,15,17,18,34,35,36,38,40,41,43,44,60,61,62,64,66,67,69,70,86,87,88,90,92,93,95,96,112,113,114,116,118,119,121,122,138,139,140,142,144,145,147,148,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)
foreach ($r in $swappedRecordedByRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# --- 3. Per-group summary table (rows 15-20): Missing/Pending counts (P, Q) ---
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 4
$ws.Range("P16").Value = 2
$ws.Range("Q16").Value = 4
$ws.Range("P17").Value = 2
$ws.Range("Q17").Value = 4
$ws.Range("P18").Value = 2
$ws.Range("Q18").Value = 4
$ws.Range("P19").Value = 2
$ws.Range("Q19").Value = 4
$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 4

# --- 4. Newly-unrecorded sessions: status flips from "Pending" (yellow) to
#         "Not Recorded" (pink) for the A:I block of each row. Copy the
#         existing "Not Recorded" formatting (style already used a couple of
#         rows above each one) onto the row, then fix up the status text.
$notRecordedRows = @(
    @{Target = 23;  Source = 21},
    @{Target = 49;  Source = 47},
    @{Target = 75;  Source = 73},
    @{Target = 101; Source = 99},
    @{Target = 127; Source = 125},
    @{Target = 153; Source = 151}
)

foreach ($pair in $notRecordedRows) {
    $t = $pair.Target
    $s = $pair.Source
    $ws.Range("A$s`:I$s").Copy()
    $ws.Range("A$t`:I$t").PasteSpecial(-4122)
    $ws.Range("I$t").Value = "Not Recorded"
}
$excel.CutCopyMode = $false
